# Insert a new weekly price record for Chirimoya (Macroferia Regional de Talca)
# as row 44, pushing the existing rows 44-54 down to rows 45-55.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 44:54 down by one to make room for the new record.
$ws.Rows(44).Insert()

# Populate the newly inserted row 44 with the new observation.
$ws.Range("A44").Value = 5
$ws.Range("B44").Value = "Macroferia Regional de Talca"
$ws.Range("C44").Value = "Maule"
$ws.Range("D44").Value = 44511
$ws.Range("E44").Value = 7
$ws.Range("F44").Value = "Fruta"
$ws.Range("G44").Value = 100107
$ws.Range("H44").Value = "Otros"
$ws.Range("I44").Value = 100107002
$ws.Range("J44").Value = "Chirimoya"
$ws.Range("K44").Value = "Cultivar IV Región"
$ws.Range("L44").Value = "Primera"
$ws.Range("M44").Value = 200
$ws.Range("N44").Value = 25000
$ws.Range("O44").Value = 25000
$ws.Range("P44").Value = 25000
$ws.Range("Q44").Value = "$/bandeja 10 kilos"
$ws.Range("R44").Value = "Provincia de Limarí"
$ws.Range("S44").Value = 2500
$ws.Range("T44").Value = 10
